$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), copying the existing
# header style (bold, bordered, centered) from H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I (I0) and J (IF), rows 2..32
$data = @(
    @(4, 6),
    @(9, 9),
    @(5, 7),
    @(6, 7),
    @(3, 5),
    @(2, 3),
    @(7, 8),
    @(10, 10),
    @(7, 7),
    @(7, 7),
    @(7, 8),
    @(9, 9),
    @(6, 7),
    @(3, 4),
    @(7, 8),
    @(8, 9),
    @(4, 5),
    @(9, 9),
    @(6, 6),
    @(1, 2),
    @(5, 6),
    @(7, 7),
    @(9, 9),
    @(6, 6),
    @(7, 8),
    @(7, 8),
    @(5, 6),
    @(2, 4),
    @(3, 4),
    @(5, 5),
    @(1, 2)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
